$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column headers to reflect 6 recall trials instead of delay columns
$ws.Range("H3").Value = "Recall 6"
$ws.Range("I3").Value = "15 Min Delay"

# Apply the same header style as B1/G1 to the new H1 cell (empty, but styled)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = $null

# Update the selected cell to match the saved selection state
$ws.Range("J8").Select() | Out-Null
